$d = $word.ActiveDocument

# 1. "increased severity however," -> "increased severity. However,"
$d.Content.Find.Execute("increased severity however, it has been shown", $false, $false, $false, $false, $false, $true, 1, $false, "increased severity. However, it has been shown", 2) | Out-Null

# 2. "utilized" -> "utilised"
$d.Content.Find.Execute("we have utilized atmospheric", $false, $false, $false, $false, $false, $true, 1, $false, "we have utilised atmospheric", 2) | Out-Null

# 3. "Nonmetirc" -> "Nonmetric"
$d.Content.Find.Execute("Nonmetirc multidimensional", $false, $false, $false, $false, $false, $true, 1, $false, "Nonmetric multidimensional", 2) | Out-Null

# 4. "during MHW different from" -> "during MHW are different from"
$d.Content.Find.Execute("during MHW different from", $false, $false, $false, $false, $false, $true, 1, $false, "during MHW are different from", 2) | Out-Null

# 5. "oceanic state, but rather through the unpredictable chaos of the climate system." -> "oceanic state that could be described by the SOM analysis."
$d.Content.Find.Execute("oceanic state, but rather through the unpredictable chaos of the climate system.", $false, $false, $false, $false, $false, $true, 1, $false, "oceanic state that could be described by the SOM analysis.", 2) | Out-Null

# 6. "wind and water patterns were shown" -> "wind and current patterns was shown"
$d.Content.Find.Execute("wind and water patterns were shown", $false, $false, $false, $false, $false, $true, 1, $false, "wind and current patterns was shown", 2) | Out-Null

# 7. Make "in situ)" italic - find the range and set italic formatting
$rng = $d.Content
$rng.Find.Execute("in situ)", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Font.Italic = $true

# 8. Set paragraph mark run properties lang (the empty <w:rPr/> -> <w:rPr><w:lang w:val="en-US"/></w:rPr>)
$para = $d.Paragraphs(1)
$para.Range.ParagraphFormat.Style = $para.Range.ParagraphFormat.Style
$markRange = $d.Range($d.Paragraphs(1).Range.End - 1, $d.Paragraphs(1).Range.End)
$markRange.LanguageID = 1033
